$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-15 16:18:37"
$ws.Range("E3").Value = "2026-02-15 16:18:39"
$ws.Range("K3").Value = "6.7 MJ/m2"
$ws.Range("O3").Value = "-5.9 °C"
$ws.Range("E4").Value = "2026-02-15 16:18:42"
$ws.Range("H4").Formula = '="69%"'
$ws.Range("H4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4163) | Out-Null
$ws.Range("K4").Value = "11.2 MJ/m2"
$ws.Range("O4").Value = "6.8 °C"
$ws.Range("E5").Value = "2026-02-15 16:18:44"
$ws.Range("K5").Value = "5.5 MJ/m2"
$ws.Range("O5").Value = "-5.3 °C"
$ws.Range("E6").Value = "2026-02-15 16:18:47"
$ws.Range("K6").Value = "11.7 MJ/m2"
$ws.Range("O6").Value = "8.1 °C"
$ws.Range("E7").Value = "2026-02-15 16:18:49"
$ws.Range("K7").Value = "11.8 MJ/m2"
$ws.Range("O7").Value = "11.3 °C"
$ws.Range("E8").Value = "2026-02-15 16:18:52"
$ws.Range("K8").Value = "11.8 MJ/m2"
$ws.Range("E9").Value = "2026-02-15 16:18:54"
$ws.Range("H9").Formula = '="45%"'
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4163) | Out-Null
$ws.Range("K9").Value = "11.6 MJ/m2"
$ws.Range("O9").Value = "11.0 °C"
$ws.Range("E10").Value = "2026-02-15 16:18:57"
$ws.Range("K10").Value = "11.7 MJ/m2"
$ws.Range("O10").Value = "7.2 °C"
$ws.Range("E11").Value = "2026-02-15 16:18:59"
$ws.Range("O11").Value = "7.5 °C"
$ws.Range("E12").Value = "2026-02-15 16:19:02"
$ws.Range("H12").Formula = '="49%"'
$ws.Range("H12").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4163) | Out-Null
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-15 16:19:04"
$ws.Range("K13").Value = "6.6 MJ/m2"
$ws.Range("O13").Value = "6.1 °C"
$ws.Range("E14").Value = "2026-02-15 16:19:07"
$ws.Range("K14").Value = "11.4 MJ/m2"
$ws.Range("O14").Value = "10.7 °C"
$ws.Range("E15").Value = "2026-02-15 16:19:09"
$ws.Range("E16").Value = "2026-02-15 16:19:12"
$ws.Range("K16").Value = "9.1 MJ/m2"
$ws.Range("O16").Value = "-2.4 °C"
$ws.Range("E17").Value = "2026-02-15 16:19:14"
$ws.Range("H17").Formula = '="32%"'
$ws.Range("H17").Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4163) | Out-Null
$ws.Range("K17").Value = "12.6 MJ/m2"
$ws.Range("E18").Value = "2026-02-15 16:19:17"
$ws.Range("K18").Value = "11.8 MJ/m2"
$ws.Range("O18").Value = "6.6 °C"
$ws.Range("E19").Value = "2026-02-15 16:19:19"
$ws.Range("K19").Value = "11.7 MJ/m2"
$ws.Range("O19").Value = "2.7 °C"
$ws.Range("E20").Value = "2026-02-15 16:19:22"
$ws.Range("H20").Formula = '="56%"'
$ws.Range("H20").Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4163) | Out-Null
$ws.Range("K20").Value = "12.8 MJ/m2"
$ws.Range("O20").Value = "-3.2 °C"
$ws.Range("E21").Value = "2026-02-15 16:19:24"
$ws.Range("J21").Value = "1015.0 hPa"
$ws.Range("K21").Value = "10.7 MJ/m2"
$ws.Range("O21").Value = "7.4 °C"
$ws.Range("E22").Value = "2026-02-15 16:19:27"
$ws.Range("K22").Value = "12.3 MJ/m2"
$ws.Range("N22").Value = "-6.5 °C 15:48 TU"
$ws.Range("E23").Value = "2026-02-15 16:19:29"
$ws.Range("K23").Value = "12.8 MJ/m2"
$ws.Range("O23").Value = "-4.2 °C"
$ws.Range("E24").Value = "2026-02-15 16:19:32"
$ws.Range("J24").Value = "1018.3 hPa"
$ws.Range("K24").Value = "11.3 MJ/m2"
$ws.Range("O24").Value = "8.2 °C"
$ws.Range("E25").Value = "2026-02-15 16:19:34"
$ws.Range("K25").Value = "8.8 MJ/m2"
$ws.Range("O25").Value = "-2.1 °C"
$ws.Range("E26").Value = "2026-02-15 16:19:37"
$ws.Range("E27").Value = "2026-02-15 16:19:39"
$ws.Range("H27").Formula = '="44%"'
$ws.Range("H27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4163) | Out-Null
$ws.Range("K27").Value = "11.2 MJ/m2"
$ws.Range("O27").Value = "-0.4 °C"
$ws.Range("E28").Value = "2026-02-15 16:19:42"
$ws.Range("K28").Value = "10.7 MJ/m2"
$ws.Range("L28").Value = "19.1 km/h - 58º 15:48 TU"
$ws.Range("O28").Value = "6.0 °C"
$ws.Range("E29").Value = "2026-02-15 16:19:44"
$ws.Range("H29").Formula = '="52%"'
$ws.Range("H29").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4163) | Out-Null
$ws.Range("K29").Value = "12.0 MJ/m2"
$ws.Range("O29").Value = "10.1 °C"
$ws.Range("E30").Value = "2026-02-15 16:19:47"
$ws.Range("K30").Value = "11.9 MJ/m2"
$ws.Range("E31").Value = "2026-02-15 16:19:49"
$ws.Range("K31").Value = "10.2 MJ/m2"
$ws.Range("O31").Value = "9.4 °C"
$ws.Range("E32").Value = "2026-02-15 16:19:52"
$ws.Range("K32").Value = "9.0 MJ/m2"
$ws.Range("O32").Value = "3.0 °C"
$ws.Range("E33").Value = "2026-02-15 16:19:54"
$ws.Range("K33").Value = "10.3 MJ/m2"
$ws.Range("O33").Value = "5.2 °C"
$ws.Range("E34").Value = "2026-02-15 16:19:56"
$ws.Range("H34").Formula = '="48%"'
$ws.Range("H34").Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4163) | Out-Null
$ws.Range("K34").Value = "12.0 MJ/m2"
$ws.Range("E35").Value = "2026-02-15 16:19:59"
$ws.Range("K35").Value = "10.2 MJ/m2"
$ws.Range("O35").Value = "3.6 °C"
$ws.Range("E36").Value = "2026-02-15 16:20:02"
$ws.Range("K36").Value = "9.9 MJ/m2"
$ws.Range("O36").Value = "11.3 °C"
$ws.Range("E37").Value = "2026-02-15 16:20:04"
$ws.Range("O37").Value = "5.9 °C"
$ws.Range("E38").Value = "2026-02-15 16:20:07"
$ws.Range("K38").Value = "12.1 MJ/m2"
$ws.Range("O38").Value = "7.2 °C"
$ws.Range("E39").Value = "2026-02-15 16:20:09"
$ws.Range("E40").Value = "2026-02-15 16:20:12"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-15 16:20:14"
$ws.Range("K41").Value = "12.2 MJ/m2"
$ws.Range("O41").Value = "11.8 °C"
$ws.Range("E42").Value = "2026-02-15 16:20:17"
$ws.Range("H42").Formula = '="51%"'
$ws.Range("H42").Copy() | Out-Null
$ws.Range("H42").PasteSpecial(-4163) | Out-Null
$ws.Range("O42").Value = "10.6 °C"
$ws.Range("E43").Value = "2026-02-15 16:20:19"
$ws.Range("H43").Formula = '="67%"'
$ws.Range("H43").Copy() | Out-Null
$ws.Range("H43").PasteSpecial(-4163) | Out-Null
$ws.Range("K43").Value = "12.6 MJ/m2"
$ws.Range("O43").Value = "5.5 °C"
$ws.Range("E44").Value = "2026-02-15 16:20:21"
$ws.Range("K44").Value = "9.4 MJ/m2"
$ws.Range("O44").Value = "-4.6 °C"
$ws.Range("E45").Value = "2026-02-15 16:20:24"
$ws.Range("J45").Value = "1024.1 hPa"
$ws.Range("K45").Value = "4.3 MJ/m2"
$ws.Range("O45").Value = "0.3 °C"
$ws.Range("E46").Value = "2026-02-15 16:20:27"
$ws.Range("J46").Value = "1019.7 hPa"
$ws.Range("K46").Value = "12.2 MJ/m2"
$ws.Range("O46").Value = "11.2 °C"

$excel.CutCopyMode = 0
Write-Output "Applied 133 changes"